$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay text (avoid Excel auto-converting numeric-looking strings to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.402.51'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '2.070.92'
$ws.Range("E3").Value = '  +3.59%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '328.16'
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").Value = '0.5182'
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("D8").Value = '0.4324'
$ws.Range("E8").Value = '  +4.39%  '
$ws.Range("D9").Value = '0.08630'
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("D10").Value = '45.97'
$ws.Range("E10").Value = '  +6.83%  '
$ws.Range("D11").Value = '1.148'
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").Value = '24.14'
$ws.Range("E12").Value = '  -2.17%  '
$ws.Range("D13").Value = '2.073.23'
$ws.Range("E13").Value = '  +3.85%  '
$ws.Range("D14").Value = '6.605'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = '7.633'
$ws.Range("E15").Value = '  +2.74%  '
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = '94.78'
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").Value = '0.00001110'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").Value = '0.06611'
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D20").Value = '18.69'
$ws.Range("E20").Value = '  -1.41%  '
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").Value = '6.202'
$ws.Range("E22").Value = '  +0.43%  '
$ws.Range("D23").Value = '30.426.94'
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("D24").Value = '12.27'
$ws.Range("E24").Value = '  +3.82%  '
$ws.Range("D25").Value = '2.303'
$ws.Range("E25").Value = '  +3.51%  '
$ws.Range("D26").Value = '2.317.75'
$ws.Range("E26").Value = '  +4.10%  '
$ws.Range("D27").Value = '22.06'
$ws.Range("E27").Value = '  -1.19%  '
$ws.Range("D28").Value = '161.11'
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("D29").Value = '2.509'
$ws.Range("E29").Value = '  +4.49%  '
$ws.Range("D30").Value = '130.35'
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("D31").Value = '1.180'
$ws.Range("E31").Value = '  +3.81%  '
$ws.Range("D32").Value = '0.1065'
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("D33").Value = '6.037'
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").Value = '3.807'
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("E35").Value = '  +12.14%  '
$ws.Range("D36").Value = '0.02547'
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("D37").Value = '9.579'
$ws.Range("E37").Value = '  +6.22%  '
$ws.Range("D38").Value = '5.436'
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").Value = '0.06570'
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("D40").Value = '12.36'
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").Value = '0.2227'
$ws.Range("E41").Value = '  +1.10%  '
$ws.Range("D42").Value = '0.6658'
$ws.Range("E42").Value = '  +0.37%  '
$ws.Range("D43").Value = '1.230'
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("D44").Value = '1.003'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.82'
$ws.Range("E45").Value = '  +1.48%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6256'
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("D47").Value = '2.181'
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("D48").Value = '3.609'
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("D49").Value = '1.225'
$ws.Range("E49").Value = '  -3.34%  '
$ws.Range("D50").Value = '81.43'
$ws.Range("E50").Value = '  +1.74%  '
$ws.Range("D51").Value = '1.175'
$ws.Range("E51").Value = '  +5.79%  '
